# TC04_Trials_Filter_Diagnosis-AdenoRectum.xlsx
# - added the Neo4j/Cypher query text into cell A2 of the "startup" sheet
#   (this is a brand new shared string, so sharedStrings.xml count/uniqueCount
#   grow and cell A2 switches from "blank, styled" to "string, styled")
# - row 2 grows tall enough (87pt) to show the multi-line query text
# - the sheet's on-open selection is moved to the query column (B2:B7)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$query = 'MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report) WITH DISTINCT c AS c, t ,a, s WHERE c.disease IN [''Adenocarcinoma of the rectum''] RETURN coalesce(c.case_id,'''') AS `Case ID` , coalesce(t.clinical_trial_designation ,'''')as `Trial Code` , coalesce(a.arm_id,'''') As `Arm` , coalesce(a.arm_drug,'''') As `Arm Treatment` , coalesce(c.disease,'''') As Diagnosis , coalesce(c.gender,'''') As Gender , coalesce(c.race,'''') As Race , coalesce(c.ethnicity,'''') As Ethnicity'

# Fill in the query used to build the companion *_Neo4jData.xlsx workbook.
# The cell already carries the wrap-text style (s="1"); setting .Value keeps
# that style and simply turns the cell into a shared-string cell.
$ws.Range("A2").Value = $query

# Grow row 2 so the long, wrapped query is readable (matches ht="87" in the
# saved worksheet XML).
$ws.Rows.Item(2).RowHeight = 87

# Leave the selection on the newly added query column when the sheet opens.
$ws.Range("B2:B7").Select()
